{"js": "// Truncate the long Java stack trace run, keeping only the first\n// two lines of the error message (through \"...is null\") and\n// dropping the NullPointerException + full stack trace that used\n// to follow it in the same run (right before the <w:br/>).\n\nconst body = context.document.body;\nbody.load(\"paragraphs\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that contains the stack trace.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"java.lang.NullPointerException\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Anchor the deletion: start right at the newline that precedes the\n  // \"java.lang.NullPointerException\" line, end right after the final\n  // stack frame line (including its trailing newline), so the\n  // paragraph break that used to follow the whole trace is preserved.\n  const startResults = target.search(\"\\njava.lang.NullPointerException\", { matchCase: true });\n  startResults.load(\"items\");\n  const endResults = target.search(\n    \"RemoteTestRunner.main(RemoteTestRunner.java:210)\\n\",\n    { matchCase: true }\n  );\n  endResults.load(\"items\");\n  await context.sync();\n\n  if (startResults.items.length > 0 && endResults.items.length > 0) {\n    const toRemove = startResults.items[0].expandTo(endResults.items[0]);\n    toRemove.insertText(\"\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Truncate the long Java stack trace run, keeping only the first\n# two lines of the error message (through \"...is null\") and\n# dropping the NullPointerException + full stack trace that used\n# to follow it in the same run (right before the line break).\n\n$d = $word.ActiveDocument\n\n$startRng = $d.Content\n$startFound = $startRng.Find.Execute(\"`njava.lang.NullPointerException\")\n\n$endRng = $d.Content\n$endFound = $endRng.Find.Execute(\"RemoteTestRunner.main(RemoteTestRunner.java:210)`n\")\n\nif ($startFound -and $endFound) {\n    $combined = $d.Range($startRng.Start, $endRng.End)\n    $combined.Text = \"\"\n}\n"}
